$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 into the two new header cells (I1, J1)
# so they pick up style index 1 (bold, bordered, centered header style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data for the two new columns (rows 2-80).
$iValues = @(7,7,8,7,8,8,6,4,9,9,5,8,9,8,5,8,7,8,4,9,9,8,8,8,3,8,8,9,7,4,7,8,8,9,6,9,4,6,8,8,7,10,6,8,6,8,7,7,7,7,6,7,8,7,5,7,7,7,4,8,7,6,7,9,7,5,7,8,7,9,10,7,4,8,5,4,4,3,3)
$jValues = @(7,7,8,7,8,8,6,5,9,9,6,8,9,8,6,8,7,8,6,9,9,8,8,8,5,8,8,9,8,5,7,8,8,9,6,9,5,6,8,8,8,10,7,8,7,8,7,7,7,7,7,7,8,8,6,7,8,7,5,8,7,6,7,9,7,5,7,8,8,9,10,7,5,8,5,4,4,3,3)

for ($k = 0; $k -lt $iValues.Count; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
